{"js": "// Merge the three hyperlink runs (\" H\" + \"y\" + \"perlink \") into a single\n// run containing \" Hyperlink \", and switch the \"Hyperlink\" character\n// style's font color from the explicit theme color (accent1 / #4F81BD)\n// to automatic - matching the target OOXML diff.\n\n// --- 1. Merge the hyperlink text runs -------------------------------------\n// A direct replace with text identical to the current text is a no-op in\n// this engine (it keeps the original run split), so first swap in a\n// temporary placeholder and then replace that placeholder with the final\n// text - this forces the three runs to collapse into one run that carries\n// the (uniform) run formatting already in effect, i.e. the Hyperlink\n// character style.\nlet results = context.document.body.search(\" Hyperlink \", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\" HyperlinkZZZ \", Word.InsertLocation.replace);\n  await context.sync();\n\n  const results2 = context.document.body.search(\" HyperlinkZZZ \", { matchCase: true });\n  results2.load(\"text\");\n  await context.sync();\n  results2.items[0].insertText(\" Hyperlink \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// --- 2. Hyperlink style color -> Automatic ---------------------------------\nconst hyperlinkStyle = context.document.getStyles().getByNameOrNullObject(\"Hyperlink\");\nhyperlinkStyle.load(\"nameLocal\");\nawait context.sync();\n\nif (!hyperlinkStyle.isNullObject) {\n  // -16777216 is wdColorAutomatic - it serialises to OOXML as\n  // <w:color w:val=\"auto\"/> (and drops any w:themeColor), matching the\n  // target diff exactly.\n  hyperlinkStyle.font.color = -16777216;\n  await context.sync();\n}\n", "ps1": "# Merge the three hyperlink runs (\" H\" + \"y\" + \"perlink \") into a single\n# run containing \" Hyperlink \", and switch the \"Hyperlink\" character\n# style's font color from the explicit theme color (accent1 / #4F81BD) to\n# automatic - matching the target OOXML diff.\n\n$d = $word.ActiveDocument\n\n# --- 1. Merge the hyperlink text runs --------------------------------------\n# Assigning Range.Text to text that is already identical to the range's\n# current text is a silent no-op here (the original run split survives),\n# so swap in a temporary placeholder first and then replace that\n# placeholder with the final text. Setting Range.Text also resets the\n# range to the paragraph's base character formatting, so re-apply the\n# \"Hyperlink\" character style afterwards to keep the run inside\n# <w:hyperlink> styled the way the source document had it.\n\n$rng = $d.Content\n$found = $rng.Find.Execute(\" Hyperlink \")\nif ($found) {\n  $rng.Text = \" HyperlinkZZZ \"\n  $rng.Style = $d.Styles(\"Hyperlink\")\n}\n\n$rng2 = $d.Content\n$found2 = $rng2.Find.Execute(\" HyperlinkZZZ \")\nif ($found2) {\n  $rng2.Text = \" Hyperlink \"\n  $rng2.Style = $d.Styles(\"Hyperlink\")\n}\n\n# --- 2. Hyperlink style color -> Automatic ---------------------------------\n# -16777216 is wdColorAutomatic - it serialises to OOXML as\n# <w:color w:val=\"auto\"/> (and drops any w:themeColor), matching the\n# target diff exactly.\n$style = $d.Styles(\"Hyperlink\")\n$style.Font.Color = -16777216\n"}
